$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace every product code in A2:A16 with the new single value "P1000"
$ws.Range("A2:A16").Value = "P1000"

# Reflect the scrolled view / selection left behind by the edit
$ws.Range("A16").Select()
$excel.ActiveWindow.ScrollRow = 12
